# Add a new "2022-Q4" sheet (with fresh fund-holding data) right after the
# "总计" (summary) sheet, pushing the existing quarterly sheets one slot to
# the right, and record the new quarter in the "总计" roll-up sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q4" sheet by duplicating "2022-Q3" (this keeps
#    all header text/styles identical) and inserting it right before the
#    sheet it was copied from, then overwrite it with the new data.
# ---------------------------------------------------------------------
$refSheet = $wb.Worksheets.Item("2022-Q3")
$refSheet.Copy($refSheet)
$newSheet = $wb.Worksheets.Item("2022-Q3 (2)")
$newSheet.Name = "2022-Q4"

# Header row - D1 label differs ("基金规模" instead of the old "基金金额")
$newSheet.Range("D1").Value = "基金规模"

# Data row - fund code / name / scale / position values must stay TEXT
# (they look numeric), so force a text number format before assignment.
$newSheet.Range("A2").Value = 0

$newSheet.Range("B2").NumberFormat = "@"
$newSheet.Range("B2").Value = "159743"

$newSheet.Range("C2").Value = "博时中证湖北新旧动能转换ETF"

$newSheet.Range("D2").NumberFormat = "@"
$newSheet.Range("D2").Value = "3.43"

$newSheet.Range("E2").NumberFormat = "@"
$newSheet.Range("E2").Value = "99.18"

$newSheet.Range("F2").NumberFormat = "@"
$newSheet.Range("F2").Value = "1.68"

$newSheet.Range("G2").NumberFormat = "@"
$newSheet.Range("G2").Value = "0.0576"

$newSheet.Range("H2").Value = 7

# ---------------------------------------------------------------------
# 2. Update the "总计" sheet: push the existing quarter rows down by one
#    and insert the 2022-Q4 summary figures at the top.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Give the brand-new last row (row 6) the same look as the other index
# cells in column A before we populate it.
$summary.Range("A5").Copy()
$summary.Range("A6").PasteSpecial(-4122)

for ($r = 5; $r -ge 2; $r--) {
    $nr = $r + 1
    $summary.Range("B$nr").Value = $summary.Range("B$r").Value2
    $summary.Range("C$nr").Value = $summary.Range("C$r").Value2
    $summary.Range("D$nr").Value = $summary.Range("D$r").Value2
}

$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 1
$summary.Range("D2").Value = 0.06

for ($r = 2; $r -le 6; $r++) {
    $summary.Range("A$r").Value = $r - 2
}

# ---------------------------------------------------------------------
# 3. Restore the originally-active tab: the last sheet ("2020-Q4") was
#    the selected tab before the edit, but creating/renaming sheets
#    above moved the "active sheet" marker onto the new sheet.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("2020-Q4").Activate()
